# Salt Lake Home Passing sheet: add a "Match ID" column at the front of the
# table (new column A), shifting the existing A:AC columns to B:AD, and
# filling the new column with the match id (18) for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column A; everything that was in A:AC
# moves to B:AD (values, shared strings, styles, merged cells all shift).
$ws.Columns("A").Insert()

# Header row (row 3 is the visible header row) gets the new column title.
$ws.Range("A3").Value = "Match ID"

# Header cell + all data rows (4-19) get the bold "label" style used
# elsewhere in the header (font bold, no border/alignment) - this mirrors
# the existing cellXfs entries and creates the new one Excel needs.
$ws.Range("A3:A19").Font.Bold = $true

# Fill in the match id for every row, including the hidden total row 20
# (which keeps the default / unbolded style).
for ($r = 4; $r -le 20; $r++) {
    $ws.Cells.Item($r, 1).Value = 18
}

# Selection moves to the newly added Match ID column's data rows.
$ws.Range("A3:A19").Select()
